$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# Fix the header row (B1:F1) which previously (erroneously) duplicated the
# first data row's values instead of containing the proper field names, and
# extend the header with the common metadata columns used by the other
# sheets in this workbook (G1:M1).
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# Populate the new metadata columns (G:M) for each existing data row (rows
# 2-7), matching the pattern used across the other worksheets in the
# workbook.
for ($row = 2; $row -le 7; $row++) {
    $idx = $ws.Cells.Item($row, 1).Value

    $ws.Cells.Item($row, 7).Value  = "deposit"
    $ws.Cells.Item($row, 8).Value  = "normal"
    $ws.Cells.Item($row, 9).Value  = "2012-04-23"
    $ws.Cells.Item($row, 10).Value = "林郁方"
    $ws.Cells.Item($row, 11).Value = 716
    $ws.Cells.Item($row, 12).Value = "tmp5c281"
    $ws.Cells.Item($row, 13).Value = $idx
}
